# Apply cell value updates per the target diff, sheet by sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7282.222
$ws.Range("J18").Value = 3180.6667
$ws.Range("L18").Value = 3180.6667
$ws.Range("N18").Value = -3748.6667
$ws.Range("H51").Value = 9988.5
$ws.Range("J51").Value = 9988
$ws.Range("L51").Value = 9988
$ws.Range("N51").Value = -10956
$ws.Range("H101").Value = 12988936
$ws.Range("I101").Value = 17858788
$ws.Range("K101").Value = 53576364
$ws.Range("M101").Value = -53574742
$ws.Range("H109").Value = 342037500
$ws.Range("J109").Value = 342037500
$ws.Range("L109").Value = 342037500
$ws.Range("N109").Value = -342040274
$ws.Range("H112").Value = 30026.428
$ws.Range("J112").Value = 39972.117
$ws.Range("L112").Value = 119916.351
$ws.Range("N112").Value = -122132.351
$ws.Range("H113").Value = 7634.7
$ws.Range("J113").Value = 9000
$ws.Range("L113").Value = 9000
$ws.Range("N113").Value = -15508
$ws.Range("H116").Value = 2227472
$ws.Range("I116").Value = 5557405
$ws.Range("K116").Value = 5557405
$ws.Range("M116").Value = -5553963
$ws.Range("H135").Value = 29998
$ws.Range("J135").Value = 9995
$ws.Range("L135").Value = 89955
$ws.Range("N135").Value = -95025
$ws.Range("H138").Value = 3461.318
$ws.Range("J138").Value = 4549.933
$ws.Range("L138").Value = 13649.799
$ws.Range("N138").Value = -23929.799
$ws.Range("H140").Value = 81123.5
$ws.Range("J140").Value = 81123.5
$ws.Range("L140").Value = 81123.5
$ws.Range("N140").Value = -91483.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1640.0741
$ws.Range("I32").Value = 1633.283
$ws.Range("K32").Value = 1633.283
$ws.Range("M32").Value = -1346.283
$ws.Range("H61").Value = 3085.1707
$ws.Range("I61").Value = 1911.7241
$ws.Range("K61").Value = 1911.7241
$ws.Range("M61").Value = -1699.7241
$ws.Range("H63").Value = 3982.3333
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 3982.3333
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 3712.9429
$ws.Range("I74").Value = 2278.88
$ws.Range("K74").Value = 2278.88
$ws.Range("M74").Value = -1404.88
$ws.Range("H77").Value = 3712.9429
$ws.Range("I77").Value = 2278.88
$ws.Range("K77").Value = 11394.4
$ws.Range("M77").Value = -7026.400000000001
$ws.Range("H136").Value = 3085.1707
$ws.Range("I136").Value = 1911.7241
$ws.Range("K136").Value = 5735.1723
$ws.Range("M136").Value = -3185.1723

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8708.454
$ws.Range("I86").Value = 14181.2
$ws.Range("K86").Value = 14181.2
$ws.Range("M86").Value = -13058.2
$ws.Range("H89").Value = 8708.454
$ws.Range("I89").Value = 14181.2
$ws.Range("K89").Value = 70906
$ws.Range("M89").Value = -65290
$ws.Range("H99").Value = 16873.576
$ws.Range("J99").Value = 12333
$ws.Range("L99").Value = 12333
$ws.Range("N99").Value = -15329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 23332.818
$ws.Range("J88").Value = 20956.777
$ws.Range("L88").Value = 20956.777
$ws.Range("N88").Value = -21768.777
$ws.Range("H91").Value = 23332.818
$ws.Range("J91").Value = 20956.777
$ws.Range("L91").Value = 20956.777
$ws.Range("N91").Value = -23764.777
$ws.Range("H93").Value = 14999
$ws.Range("I93").Value = 14999
$ws.Range("K93").Value = 14999
$ws.Range("M93").Value = -13127

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 37493852
$ws.Range("I4").Value = 29163166
$ws.Range("K4").Value = 87489498
$ws.Range("M4").Value = -87489386
$ws.Range("H48").Value = 5721.636
$ws.Range("I48").Value = 1451.2
$ws.Range("J48").Value = 9280.333000000001
$ws.Range("K48").Value = 4353.6
$ws.Range("L48").Value = 27840.999
$ws.Range("M48").Value = -4103.6
$ws.Range("N48").Value = -28340.999
$ws.Range("H49").Value = 699.5
$ws.Range("I49").Value = 699.5
$ws.Range("K49").Value = 2098.5
$ws.Range("M49").Value = -1942.5
$ws.Range("H114").Value = 3999.4285
$ws.Range("J114").Value = 4249.5
$ws.Range("L114").Value = 12748.5
$ws.Range("N114").Value = -19256.5
$ws.Range("H121").Value = 5000325
$ws.Range("I121").Value = 433
$ws.Range("K121").Value = 1299
$ws.Range("M121").Value = 11

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6075.2593
$ws.Range("I70").Value = 5573.5
$ws.Range("K70").Value = 5573.5
$ws.Range("M70").Value = -5303.5
$ws.Range("H73").Value = 6075.2593
$ws.Range("I73").Value = 5573.5
$ws.Range("K73").Value = 5573.5
$ws.Range("M73").Value = -4637.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 8004.5
$ws.Range("I4").Value = 6009
$ws.Range("K4").Value = 6009
$ws.Range("M4").Value = -5896
$ws.Range("H20").Value = 6502.5
$ws.Range("I20").Value = 7005
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 7005
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -6779
$ws.Range("N20").Value = -6452
$ws.Range("H22").Value = 3423.8572
$ws.Range("I22").Value = 3826
$ws.Range("K22").Value = 3826
$ws.Range("M22").Value = -3531
$ws.Range("H27").Value = 3423.8572
$ws.Range("I27").Value = 3826
$ws.Range("K27").Value = 3826
$ws.Range("M27").Value = -3719
$ws.Range("H28").Value = 8004.5
$ws.Range("I28").Value = 6009
$ws.Range("K28").Value = 6009
$ws.Range("M28").Value = -5777
$ws.Range("H37").Value = 8004.5
$ws.Range("I37").Value = 6009
$ws.Range("K37").Value = 6009
$ws.Range("M37").Value = -5902
$ws.Range("H46").Value = 3767.6365
$ws.Range("I46").Value = 908.1667
$ws.Range("J46").Value = 7199
$ws.Range("K46").Value = 908.1667
$ws.Range("L46").Value = 7199
$ws.Range("M46").Value = -720.1667
$ws.Range("N46").Value = -7575
$ws.Range("H96").Value = 36666.332
$ws.Range("J96").Value = 36666.332
$ws.Range("L96").Value = 36666.332
$ws.Range("N96").Value = -42158.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2165
$ws.Range("I14").Value = 2997.5
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 2997.5
$ws.Range("L14").Value = 500
$ws.Range("M14").Value = -2829.5
$ws.Range("N14").Value = -836
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H20").Value = 10755
$ws.Range("J20").Value = 1500
$ws.Range("L20").Value = 1500
$ws.Range("N20").Value = -1980
$ws.Range("H32").Value = 16508.25
$ws.Range("I32").Value = 14581
$ws.Range("J32").Value = 29999
$ws.Range("K32").Value = 14581
$ws.Range("L32").Value = 29999
$ws.Range("M32").Value = -14264
$ws.Range("N32").Value = -30633
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H124").Value = 47998.332
$ws.Range("J124").Value = 47998.332
$ws.Range("L124").Value = 47998.332
$ws.Range("N124").Value = -57818.332
$ws.Range("H136").Value = 2312.6316
$ws.Range("I136").Value = 1746.0667
$ws.Range("J136").Value = 4437.25
$ws.Range("K136").Value = 5238.2001
$ws.Range("L136").Value = 13311.75
$ws.Range("M136").Value = -2688.2001
$ws.Range("N136").Value = -18411.75
